# Applies the author's latest edits to threshold.xlsx:
#  - updates the "cb_distances_calc" max (C3) from 6.3 to 6
#  - updates the "ratio" max (C4) from 27 to 25
#  - leaves the active selection on J13 (the cell selected when the
#    workbook was last saved)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C3").Value = 6
$ws.Range("C4").Value = 25

$ws.Range("J13").Select()
